$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Header label renames on the survey sheet (order matters for shared-string
# table append order: hint.text must be inserted before prompt.text)
$survey.Range("G1").Value = "display.hint.text"
$survey.Range("F1").Value = "display.prompt.text"

# Settings sheet label rename
$settings.Range("C1").Value = "display.title.text"

# Selections / active sheet changed: settings sheet becomes active/selected
[void]$survey.Range("F2").Select()
[void]$settings.Range("C10").Select()
[void]$settings.Select()
